$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 7818668
$ws.Range("I62").Value = 13893632
$ws.Range("K62").Value = 13893632
$ws.Range("M62").Value = -13893008
$ws.Range("H65").Value = 7818668
$ws.Range("I65").Value = 13893632
$ws.Range("K65").Value = 69468160
$ws.Range("M65").Value = -69465040
$ws.Range("H132").Value = 2316.35
$ws.Range("I132").Value = 2309.5676
$ws.Range("J132").Value = 2400
$ws.Range("K132").Value = 6928.702799999999
$ws.Range("L132").Value = 7200
$ws.Range("M132").Value = -4398.702799999999
$ws.Range("N132").Value = -12260
$ws.Range("H137").Value = 2736.6191
$ws.Range("I137").Value = 1981.9286
$ws.Range("J137").Value = 4246
$ws.Range("K137").Value = 5945.7858
$ws.Range("L137").Value = 12738
$ws.Range("M137").Value = -3395.7858
$ws.Range("N137").Value = -17838

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3999.2
$ws.Range("I32").Value = 3362.8333
$ws.Range("K32").Value = 3362.8333
$ws.Range("M32").Value = -3075.8333
$ws.Range("H61").Value = 4532.15
$ws.Range("I61").Value = 4586.839
$ws.Range("J61").Value = 4343.778
$ws.Range("K61").Value = 4586.839
$ws.Range("L61").Value = 4343.778
$ws.Range("M61").Value = -4374.839
$ws.Range("N61").Value = -4767.778
$ws.Range("H63").Value = 5682.905
$ws.Range("I63").Value = 2861.25
$ws.Range("J63").Value = 9445.111000000001
$ws.Range("K63").Value = 2861.25
$ws.Range("L63").Value = 9445.111000000001
$ws.Range("M63").Value = -2175.25
$ws.Range("N63").Value = -10817.111
$ws.Range("H66").Value = 5682.905
$ws.Range("I66").Value = 2861.25
$ws.Range("J66").Value = 9445.111000000001
$ws.Range("K66").Value = 14306.25
$ws.Range("L66").Value = 47225.55500000001
$ws.Range("M66").Value = -10874.25
$ws.Range("N66").Value = -54089.55500000001
$ws.Range("H74").Value = 3131.1
$ws.Range("I74").Value = 3740.6155
$ws.Range("K74").Value = 3740.6155
$ws.Range("M74").Value = -2866.6155
$ws.Range("H77").Value = 3131.1
$ws.Range("I77").Value = 3740.6155
$ws.Range("K77").Value = 18703.0775
$ws.Range("M77").Value = -14335.0775
$ws.Range("H110").Value = 209712.92
$ws.Range("I110").Value = 218730.88
$ws.Range("K110").Value = 218730.88
$ws.Range("M110").Value = -216685.88
$ws.Range("H122").Value = 4685.6
$ws.Range("I122").Value = 3220.0833
$ws.Range("K122").Value = 9660.249899999999
$ws.Range("M122").Value = -7210.249899999999
$ws.Range("H132").Value = 4434.328
$ws.Range("I132").Value = 3393.238
$ws.Range("K132").Value = 10179.714
$ws.Range("M132").Value = -7649.714
$ws.Range("H136").Value = 4532.15
$ws.Range("I136").Value = 4586.839
$ws.Range("J136").Value = 4343.778
$ws.Range("K136").Value = 13760.517
$ws.Range("L136").Value = 13031.334
$ws.Range("M136").Value = -11210.517
$ws.Range("N136").Value = -18131.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 21279262
$ws.Range("I20").Value = 27780292
$ws.Range("K20").Value = 27780292
$ws.Range("M20").Value = -27780045
$ws.Range("H86").Value = 775237.25
$ws.Range("I86").Value = 1309524.8
$ws.Range("J86").Value = 3488.5557
$ws.Range("K86").Value = 1309524.8
$ws.Range("L86").Value = 3488.5557
$ws.Range("M86").Value = -1308401.8
$ws.Range("N86").Value = -5734.5557
$ws.Range("H89").Value = 775237.25
$ws.Range("I89").Value = 1309524.8
$ws.Range("J89").Value = 3488.5557
$ws.Range("K89").Value = 6547624
$ws.Range("L89").Value = 17442.7785
$ws.Range("M89").Value = -6542008
$ws.Range("N89").Value = -28674.7785
$ws.Range("H105").Value = 43679.457
$ws.Range("I105").Value = 47467.637
$ws.Range("K105").Value = 47467.637
$ws.Range("M105").Value = -45720.637
$ws.Range("H134").Value = 32594.027
$ws.Range("I134").Value = 5237.2383
$ws.Range("J134").Value = 68499.81
$ws.Range("K134").Value = 15711.7149
$ws.Range("L134").Value = 205499.43
$ws.Range("M134").Value = -13176.7149
$ws.Range("N134").Value = -210569.43

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4556.0557
$ws.Range("I99").Value = 3692
$ws.Range("K99").Value = 3692
$ws.Range("M99").Value = -2194
$ws.Range("H105").Value = 711.7143
$ws.Range("I105").Value = 665.75
$ws.Range("K105").Value = 665.75
$ws.Range("M105").Value = 1081.25
$ws.Range("H107").Value = 584.3077
$ws.Range("I107").Value = 565.6
$ws.Range("J107").Value = 646.6667
$ws.Range("K107").Value = 565.6
$ws.Range("L107").Value = 646.6667
$ws.Range("M107").Value = 1354.4
$ws.Range("N107").Value = -4486.6667
$ws.Range("H122").Value = 2592.8333
$ws.Range("I122").Value = 1897.25
$ws.Range("K122").Value = 5691.75
$ws.Range("M122").Value = -3241.75
$ws.Range("H126").Value = 4556.0557
$ws.Range("I126").Value = 3692
$ws.Range("K126").Value = 11076
$ws.Range("M126").Value = -8606
$ws.Range("H132").Value = 2464.5264
$ws.Range("I132").Value = 2119.9333
$ws.Range("K132").Value = 6359.7999
$ws.Range("M132").Value = -3829.7999
$ws.Range("H134").Value = 479621.56
$ws.Range("I134").Value = 3602.65
$ws.Range("J134").Value = 10000000
$ws.Range("K134").Value = 10807.95
$ws.Range("L134").Value = 30000000
$ws.Range("M134").Value = -8272.950000000001
$ws.Range("N134").Value = -30005070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 116.21429
$ws.Range("I14").Value = 116.21429
$ws.Range("K14").Value = 348.64287
$ws.Range("M14").Value = -175.64287
$ws.Range("H56").Value = 6457.1665
$ws.Range("I56").Value = 6457.1665
$ws.Range("K56").Value = 6457.1665
$ws.Range("M56").Value = -5927.1665
$ws.Range("H127").Value = 1375.7
$ws.Range("J127").Value = 1375.7
$ws.Range("L127").Value = 4127.1
$ws.Range("N127").Value = -14047.1
$ws.Range("H132").Value = 427283.34
$ws.Range("I132").Value = 113694.664
$ws.Range("J132").Value = 593300.9
$ws.Range("K132").Value = 1023251.976
$ws.Range("L132").Value = 5339708.100000001
$ws.Range("M132").Value = -1020721.976
$ws.Range("N132").Value = -5344768.100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 262.94736
$ws.Range("I2").Value = 249.77777
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 249.77777
$ws.Range("L2").Value = 500
$ws.Range("M2").Value = -136.77777
$ws.Range("N2").Value = -726
$ws.Range("H70").Value = 23752.375
$ws.Range("J70").Value = 32503
$ws.Range("L70").Value = 32503
$ws.Range("N70").Value = -33043
$ws.Range("H73").Value = 23752.375
$ws.Range("J73").Value = 32503
$ws.Range("L73").Value = 32503
$ws.Range("N73").Value = -34375
$ws.Range("H80").Value = 574379.6
$ws.Range("I80").Value = 529319.0600000001
$ws.Range("K80").Value = 529319.0600000001
$ws.Range("M80").Value = -528321.0600000001
$ws.Range("H83").Value = 574379.6
$ws.Range("I83").Value = 529319.0600000001
$ws.Range("K83").Value = 2646595.3
$ws.Range("M83").Value = -2641603.3
$ws.Range("H122").Value = 3444.68
$ws.Range("I122").Value = 2407.0588
$ws.Range("J122").Value = 5649.625
$ws.Range("K122").Value = 7221.176399999999
$ws.Range("L122").Value = 16948.875
$ws.Range("M122").Value = -4771.176399999999
$ws.Range("N122").Value = -21848.875
$ws.Range("H126").Value = 3946.5715
$ws.Range("I126").Value = 3469.25
$ws.Range("J126").Value = 4137.5
$ws.Range("K126").Value = 10407.75
$ws.Range("L126").Value = 12412.5
$ws.Range("M126").Value = -7937.75
$ws.Range("N126").Value = -17352.5
$ws.Range("H132").Value = 33007.027
$ws.Range("I132").Value = 6785.08
$ws.Range("K132").Value = 20355.24
$ws.Range("M132").Value = -17825.24

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1241.75
$ws.Range("I55").Value = 306.83334
$ws.Range("J55").Value = 1802.7
$ws.Range("K55").Value = 306.83334
$ws.Range("L55").Value = 1802.7
$ws.Range("M55").Value = -133.83334
$ws.Range("N55").Value = -2148.7
$ws.Range("H108").Value = 54832.832
$ws.Range("J108").Value = 54832.832
$ws.Range("L108").Value = 54832.832
$ws.Range("N108").Value = -62512.832
$ws.Range("H132").Value = 3276.7896
$ws.Range("I132").Value = 2566.6428
$ws.Range("K132").Value = 7699.928400000001
$ws.Range("M132").Value = -5169.928400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 36500
$ws.Range("J70").Value = 36500
$ws.Range("L70").Value = 36500
$ws.Range("N70").Value = -37130
$ws.Range("H73").Value = 36500
$ws.Range("J73").Value = 36500
$ws.Range("L73").Value = 36500
$ws.Range("N73").Value = -38684
$ws.Range("H107").Value = 1505.8235
$ws.Range("I107").Value = 1776.9231
$ws.Range("J107").Value = 624.75
$ws.Range("K107").Value = 5330.7693
$ws.Range("L107").Value = 1874.25
$ws.Range("M107").Value = -3410.7693
$ws.Range("N107").Value = -5714.25
$ws.Range("H122").Value = 21280460
$ws.Range("I122").Value = 28575184
$ws.Range("K122").Value = 85725552
$ws.Range("M122").Value = -85723102
$ws.Range("H126").Value = 1140.6666
$ws.Range("I126").Value = 1192.1666
$ws.Range("J126").Value = 831.6667
$ws.Range("K126").Value = 3576.4998
$ws.Range("L126").Value = 2495.0001
$ws.Range("M126").Value = -1106.4998
$ws.Range("N126").Value = -7435.0001
$ws.Range("H132").Value = 2654.8708
$ws.Range("I132").Value = 2230.5952
$ws.Range("J132").Value = 3545.85
$ws.Range("K132").Value = 6691.785600000001
$ws.Range("L132").Value = 10637.55
$ws.Range("M132").Value = -4161.785600000001
$ws.Range("N132").Value = -15697.55
$ws.Range("H135").Value = 142907860
$ws.Range("J135").Value = 142907860
$ws.Range("L135").Value = 142907860
$ws.Range("N135").Value = -142918000
